$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '289.21'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.93%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '29.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.06%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.129'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '4.16%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06684'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.10%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.348'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.40%'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.356'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.71%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9174'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.35%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1586'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.55%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06674'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.13%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07666'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.02%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02964'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.54%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08995'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.11%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001592'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.17%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04504'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.52%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'One'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006475'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.99%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006290'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '2.98%'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.453'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.92%'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.406'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.27%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.220'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.97%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-2.94%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.060'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.42%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1549'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.97%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001190'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.47%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004133'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.63%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001248'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '5.77%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-1.18%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.85%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006713'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-3.77%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1244'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-12.15%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.001977'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.97%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01170'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-6.25%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005613'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.32%'
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'BOLO'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.968'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '25.94%'
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01305'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-29.43%'
